$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits
#    right after the title heading "Play Candy Bars Slot Free – Review
#    & Rating".
# -----------------------------------------------------------------
$metaPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -match "^Meta description") {
        $metaPara = $d.Paragraphs($i)
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete()
}

# -----------------------------------------------------------------
# 2) Insert a new paragraph right before the final paragraph (the one
#    that holds the italic "image prompt" text) containing a bold
#    "Play Candy Bars Slot Free – Review & Rating" run.
# -----------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$lastPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs($count)
$newPara.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Candy Bars Slot Free – Review &amp; Rating</w:t></w:r></w:p>')

# -----------------------------------------------------------------
# 3) Replace the image-prompt text in the (now last) paragraph with the
#    meta-description text, keeping its italic run formatting intact.
# -----------------------------------------------------------------
$count = $d.Paragraphs.Count
$finalPara = $d.Paragraphs($count)
$oldText = "Create a vibrant feature image for Candy Bars that features a happy Maya warrior wearing glasses, in a cartoon style. The background should be bright and colorful, with a candy-themed design such as candy canes, gumdrops, and lollipops. The Maya warrior should be holding a big lollipop and have a big smile on their face, with candy symbols surrounding them such as gumballs and chocolate bars. The image should showcase the fun and playful nature of the game while incorporating its candy theme and the idea of winning big."
$newText = "Read our review of Candy Bars slot game and play it for free. Exciting progressive jackpots and Wild symbols with multipliers. Classic slot game feel."
$finalPara.Range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
